# Updated cryptos list (price/volume refresh) per commit "Updated cryptos list on
# Tue Nov 14 22:53:59 UTC 2023 with GitHub Actions". Column D = Price, column E =
# Volume(1h), both stored as plain text in the sheet (not numbers), so numeric-
# looking price strings are written with a leading apostrophe to keep Excel from
# coercing them into numbers (which would also strip meaningful trailing zeros,
# e.g. "60.20" -> 60.2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "35.584.91"

$ws.Cells.Item(3, 4).Value = "1.986.95"
$ws.Cells.Item(3, 5).Value = "  -3.77%  "

$ws.Cells.Item(4, 5).Value = "  -0.11%  "

$ws.Cells.Item(5, 4).Value = "'242.32"
$ws.Cells.Item(5, 5).Value = "  +0.17%  "

$ws.Cells.Item(6, 5).Value = "  -3.55%  "

$ws.Cells.Item(7, 4).Value = "'57.16"
$ws.Cells.Item(7, 5).Value = "  +8.37%  "

$ws.Cells.Item(8, 5).Value = "  -0.01%  "

$ws.Cells.Item(9, 4).Value = "'60.20"
$ws.Cells.Item(9, 5).Value = "  +2.12%  "

$ws.Cells.Item(10, 4).Value = "'0.359"
$ws.Cells.Item(10, 5).Value = "  +0.01%  "

$ws.Cells.Item(11, 4).Value = "'0.0729"
$ws.Cells.Item(11, 5).Value = "  -2.91%  "

$ws.Cells.Item(12, 4).Value = "'0.102"
$ws.Cells.Item(12, 5).Value = "  -5.44%  "

$ws.Cells.Item(13, 4).Value = "'0.920"
$ws.Cells.Item(13, 5).Value = "  +2.36%  "

$ws.Cells.Item(14, 5).Value = "  -3.51%  "

$ws.Cells.Item(15, 4).Value = "2.275.03"
$ws.Cells.Item(15, 5).Value = "  -3.77%  "

$ws.Cells.Item(16, 4).Value = "'5.22"
$ws.Cells.Item(16, 5).Value = "  -3.22%  "

$ws.Cells.Item(17, 4).Value = "1.985.46"
$ws.Cells.Item(17, 5).Value = "  -3.73%  "

$ws.Cells.Item(18, 4).Value = "'17.18"
$ws.Cells.Item(18, 5).Value = "  +5.31%  "

$ws.Cells.Item(19, 4).Value = "35.522.17"
$ws.Cells.Item(19, 5).Value = "  -2.64%  "

$ws.Cells.Item(20, 5).Value = "  -2.05%  "

$ws.Cells.Item(21, 4).Value = "0.0₃0836"
$ws.Cells.Item(21, 5).Value = "  -2.99%  "

$ws.Cells.Item(22, 4).Value = "'232.81"
$ws.Cells.Item(22, 5).Value = "  -1.67%  "

$ws.Cells.Item(23, 5).Value = "  -3.62%  "

$ws.Cells.Item(24, 5).Value = "  +0.11%  "

$ws.Cells.Item(25, 5).Value = "  -2.32%  "

$ws.Cells.Item(26, 4).Value = "'2.33"
$ws.Cells.Item(26, 5).Value = "  +9.31%  "

$ws.Cells.Item(27, 4).Value = "'163.64"
$ws.Cells.Item(27, 5).Value = "  -0.18%  "

$ws.Cells.Item(28, 4).Value = "'9.11"
$ws.Cells.Item(28, 5).Value = "  -3.14%  "

$ws.Cells.Item(29, 4).Value = "'19.51"
$ws.Cells.Item(29, 5).Value = "  -4.70%  "

$ws.Cells.Item(30, 5).Value = "  -2.52%  "

$ws.Cells.Item(31, 5).Value = "  -0.90%  "

$ws.Cells.Item(32, 5).Value = "  -4.64%  "

$ws.Cells.Item(33, 4).Value = "'0.0587"
$ws.Cells.Item(33, 5).Value = "  -1.55%  "

$ws.Cells.Item(34, 5).Value = "  +10.04%  "

$ws.Cells.Item(35, 5).Value = "  -6.88%  "

$ws.Cells.Item(36, 5).Value = "  +1.94%  "

$ws.Cells.Item(37, 5).Value = "  -0.04%  "

$ws.Cells.Item(38, 5).Value = "  -2.00%  "

$ws.Cells.Item(39, 4).Value = "'4.90"
$ws.Cells.Item(39, 5).Value = "  +1.04%  "

$ws.Cells.Item(40, 5).Value = "  -4.98%  "

$ws.Cells.Item(41, 5).Value = "  -3.29%  "

$ws.Cells.Item(42, 5).Value = "  -2.76%  "

$ws.Cells.Item(43, 5).Value = "  -4.10%  "

# Row 44: rows 44/45 swapped (Aave <-> Cronos)
$ws.Cells.Item(44, 2).Value = "Cronos"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(44, 4).Value = "'0.0890"
$ws.Cells.Item(44, 5).Value = "  -4.78%  "

# Row 45: rows 44/45 swapped (Aave <-> Cronos)
$ws.Cells.Item(45, 2).Value = "Aave"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(45, 4).Value = "'90.98"
$ws.Cells.Item(45, 5).Value = "  -3.34%  "

$ws.Cells.Item(46, 4).Value = "1.377.06"
$ws.Cells.Item(46, 5).Value = "  -1.33%  "

$ws.Cells.Item(47, 4).Value = "'7.35"
$ws.Cells.Item(47, 5).Value = "  -0.60%  "

$ws.Cells.Item(48, 4).Value = "'15.44"
$ws.Cells.Item(48, 5).Value = "  -0.77%  "

$ws.Cells.Item(49, 4).Value = "'2.88"
$ws.Cells.Item(49, 5).Value = "  +0.93%  "

$ws.Cells.Item(50, 5).Value = "  -3.59%  "

$ws.Cells.Item(51, 4).Value = "'45.74"
$ws.Cells.Item(51, 5).Value = "  +1.72%  "
